$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.154.19'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.538.03'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''617.25'
$ws.Range('E5').Value = '  +5.72%  '
$ws.Range('D6').Value = '''186.01'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.217'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').Value = '''0.657'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('D11').Value = '''53.53'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E12').Value = '  -3.42%  '
$ws.Range('D13').Value = '''9.60'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').Value = '4.099.42'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '''618.18'
$ws.Range('E15').Value = '  +8.80%  '
$ws.Range('D16').Value = '70.194.75'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '''12.82'
$ws.Range('E17').Value = '  +3.45%  '
$ws.Range('D18').Value = '''18.99'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').Value = '3.525.21'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = '''0.996'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').Value = '''17.53'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('D23').Value = '''103.51'
$ws.Range('E23').Value = '  +9.20%  '
$ws.Range('D24').Value = '''4.71'
$ws.Range('E24').Value = '  +2.58%  '
$ws.Range('D25').Value = '''5.04'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('E26').Value = '  +3.64%  '
$ws.Range('D27').Value = '''10.98'
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('D28').Value = '''9.82'
$ws.Range('E28').Value = '  +8.10%  '
$ws.Range('D29').Value = '''33.87'
$ws.Range('E29').Value = '  +5.14%  '
$ws.Range('E30').Value = '  -2.83%  '
$ws.Range('D31').Value = '''12.39'
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('D33').Value = '''64.18'
$ws.Range('D34').Value = '''3.60'
$ws.Range('E34').Value = '  +16.86%  '
$ws.Range('D35').Value = '''3.21'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('D36').Value = '''532.47'
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('E37').Value = '  -2.89%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = '''37.29'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').Value = '''3.58'
$ws.Range('E40').Value = '  +6.24%  '
$ws.Range('E41').Value = '  -2.98%  '
$ws.Range('D42').Value = '3.532.92'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('E43').Value = '  +0.77%  '
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('E46').Value = '  +3.93%  '
$ws.Range('D47').Value = '''3.37'
$ws.Range('E47').Value = '  -4.92%  '
$ws.Range('D48').Value = '''9.06'
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('D51').Value = '''133.95'
$ws.Range('E51').Value = '  -1.16%  '
